# Word COM-interop script implementing the "Quantum Entanglement" -> "Chemistry"
# content rewrite described by the unified diff.
#
# Strategy: use Find/Replace (Range.Find.Execute) for nearly every text change,
# since it correctly finds text that spans run boundaries and writes back a
# single run with the first matched run's formatting (mirroring what Word
# itself does when you retype/replace a sentence). Manual line breaks in the
# replacement text use the "^l" Find/Replace code. The one sentence that
# contains a literal straight apostrophe is inserted with Range.InsertAfter
# instead, because Find/Replace silently "smart-quotes" apostrophes in its
# ReplaceWith argument.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------------
Replace-Text "Unraveling the Enigmatic Beauty: Quantum Entanglement" `
             "The Enigmatic World of Chemistry: Unveiling the Secrets of Matter"

# ---------------------------------------------------------------------------
# Byline ("Dr. Chloe Anderson" -> "Olivia Curtis")
# ---------------------------------------------------------------------------
Replace-Text "Dr. Chloe Anderson" "Olivia Curtis"

# ---------------------------------------------------------------------------
# Email line
# ---------------------------------------------------------------------------
Replace-Text "chloe" "oliviacurtis"
Replace-Text "anderson@quantumresearch" "chemistry@validmail"

# ---------------------------------------------------------------------------
# Introduction paragraph
# ---------------------------------------------------------------------------
Replace-Text "In the realm of quantum physics, where uncertainty reigns and particles defy classical logic, lies a captivating phenomenon known as quantum entanglement" `
             "In the vast expanse of the universe, chemistry stands as a pillar of understanding, unveiling the fundamental principles that govern the intricate interactions of matter at its most basic level"

Replace-Text 'This extraordinary phenomenon, often labeled as "spooky action at a distance," by Albert Einstein, challenges our intuitive understanding of reality and opens up a new realm of scientific exploration' `
             'It embarks upon an exploration of the diverse elements and compounds that make up our world, deciphering the enigmatic language of chemical reactions and revealing the profound impact they have on our existence. Chemistry weaves its way through fabrics of our lives, touching every aspect from our clothes, and medicine to our food and technology, unlocking the secrets of the molecular dance that shapes our reality'

Replace-Text "Quantum entanglement is a captivating dance between two or more particles whose properties, such as spin, polarization, or energy, become correlated in an inexplicable manner" `
             "Introduction Continued:^l^lThe history of chemistry is a rich tapestry woven with tales of intrepid pioneers who dared to question the enigmatic nature of matter"

Replace-Text "These particles remain intimately connected, regardless of the distance separating them, sharing their fates in a profound and nonlocal way" `
             "PLACEHOLDER_ALCHEMISTS_SENTENCE"

Replace-Text "The measurement of the property of one entangled particle instantaneously influences the state of the other, irrespective of the vast cosmic expanse separating them" `
             "Each discovery, each breakthrough, has added a brushstroke of color to the vibrant canvas of our understanding, revealing the hidden harmonies of the chemical world"

Replace-Text "Einstein famously declared this phenomenon as `"spooky action at a distance,`" highlighting its ethereal and counterintuitive nature" `
             "Introduction Continued:^l^lChemistry is not merely a body of knowledge; it is an art form, a symphony of elements and compounds that harmonize in countless ways"

Replace-Text "Quantum entanglement has ignited intense scientific debates, challenging our fundamental understanding of physics and pushing the boundaries of human knowledge" `
             "It is the molecular dance of atoms, the delicate interplay of forces that shape the world around us. Chemistry teaches us to observe, to experiment, to analyze, to create. It cultivates a sense of wonder and curiosity, encouraging us to ask questions about the world around us and to seek answers in the intricate web of chemical processes that underlie all of nature"

# The "philosopher's stone" sentence has a literal straight apostrophe, which
# Find/Replace would otherwise convert to a curly quote -- insert it directly
# via Range.InsertAfter (positioned via Find, without using Replace) instead.
$rng = $d.Content
$rng.Find.Execute("PLACEHOLDER_ALCHEMISTS_SENTENCE") | Out-Null
$rng.Text = ""
$ins = $rng.Duplicate
$ins.Collapse(0)
$ins.InsertAfter("From ancient alchemists seeking the elusive philosopher's stone to modern chemists unraveling the intricate secrets of DNA, the quest for knowledge has driven the evolution of this field")

# ---------------------------------------------------------------------------
# Summary heading gets a lastRenderedPageBreak marker ahead of its text.
# ---------------------------------------------------------------------------
$headingRng = $d.Content
$headingRng.Find.Execute("Summary") | Out-Null
$headingRng.Collapse(1)
$headingRng.InsertBefore([char]12)

# ---------------------------------------------------------------------------
# Summary body paragraph
# ---------------------------------------------------------------------------
Replace-Text "Quantum entanglement stands as an enigmatic masterpiece of nature, a symphony of subatomic choreography that weaves an intricate tapestry of interconnectedness" `
             "Chemistry delves into the enigmatic secrets of matter, unveiling the profound impact that chemical reactions have on our lives"

Replace-Text "It redefines our notions of locality and causality, ushering in a profound transformation in how we perceive the universe" `
             "It is a tapestry woven with the threads of history, where the quest for knowledge has driven the evolution of this field"

Replace-Text "While its underlying mechanisms may still elude our grasp, quantum entanglement holds the promise of unlocking transformative technologies and reshaping our understanding of the quantum realm. From cryptography to computation, the implications of quantum entanglement reverberate across diverse fields, beckoning us to explore this extraordinary phenomenon and its captivating implications for science, technology, and human knowledge" `
             "Chemistry is not merely a body of knowledge, but an art form, a symphony of elements and compounds that harmonize in countless ways, revealing the hidden harmonies of the chemical world"

# ---------------------------------------------------------------------------
# New trailing empty paragraph at the end of the document body.
# ---------------------------------------------------------------------------
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null

Write-Output "edit complete"
